$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "price" column header in D1
$ws.Range("D1").Value = "price"

# Fill D2:D12 with quantity (column C) * 10
for ($row = 2; $row -le 12; $row++) {
    $qty = $ws.Cells.Item($row, 3).Value2
    $ws.Cells.Item($row, 4).Value = $qty * 10
}

# Update the selection to match the recorded end-state
$ws.Range("H7").Select()
